$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mislabeled column names in A22 and A24
$ws.Range("A22").Value = "Y_COORD_CD "
$ws.Range("A24").Value = "Longitude"

# Update the active selection to B28 (as recorded in the saved view state)
$ws.Range("B28").Select()
